$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 13 ("Safety Check Inc.") - subsequent rows shift up
$ws.Rows("13:13").Delete()

# Select the row that now occupies row 13 (matches author's recorded selection)
$ws.Range("A13:XFD13").Select()
